$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 245
$ws1.Range("F3").Value = 1449
$ws1.Range("F4").Value = 22
$ws1.Range("F5").Value = 898
$ws1.Range("F7").Value = 1268
$ws1.Range("F8").Value = 1620
$ws1.Range("F11").Value = 2301
$ws1.Range("F12").Value = 465
$ws1.Range("F13").Value = 132
$ws1.Range("F15").Value = 29
$ws1.Range("F17").Value = 88
$ws1.Range("F18").Value = 6329
$ws1.Range("F19").Value = 48
$ws1.Range("F20").Value = 6268
$ws1.Range("F21").Value = 10245
$ws1.Range("F22").Value = 122
$ws1.Range("F24").Value = 187
$ws1.Range("F25").Value = 281
$ws1.Range("F26").Value = 510
$ws1.Range("F29").Value = 4406
$ws1.Range("F30").Value = 161
$ws1.Range("F31").Value = 397

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 1162

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 245
$ws4.Range("F6").Value = 1449
$ws4.Range("F7").Value = 1449
$ws4.Range("F8").Value = 22
$ws4.Range("F10").Value = 898
$ws4.Range("F12").Value = 1268
$ws4.Range("F14").Value = 1620
$ws4.Range("F17").Value = 2301
$ws4.Range("F19").Value = 465
$ws4.Range("F20").Value = 132
$ws4.Range("F22").Value = 29
$ws4.Range("F25").Value = 88
$ws4.Range("F26").Value = 6329
$ws4.Range("F27").Value = 48
$ws4.Range("F28").Value = 6268
$ws4.Range("F29").Value = 10246
$ws4.Range("F31").Value = 122
$ws4.Range("F33").Value = 187
$ws4.Range("F34").Value = 281
$ws4.Range("F36").Value = 510
$ws4.Range("F41").Value = 4406
$ws4.Range("F43").Value = 161
$ws4.Range("F48").Value = 397
